$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": mark the da49f725 and df1634af rows as handed off
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-07 06:38:15"

$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-07 06:38:15"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("H4").Value = "2016-09-07 06:38:03"
$wsZhCn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a054dbbfb576e41805903624e40ef6b845f30dbd/e2e/da49f725-a87d-4da0-87bd-44b5adf3ff8e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3ffcc59c70c94baa481c8c9b824af0f0aa261ee/e2e/da49f725-a87d-4da0-87bd-44b5adf3ff8e.md."

$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("H5").Value = "2016-09-07 06:38:03"
$wsZhCn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a054dbbfb576e41805903624e40ef6b845f30dbd/e2e/df1634af-b494-4138-8392-4b645d3171bd.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3ffcc59c70c94baa481c8c9b824af0f0aa261ee/e2e/df1634af-b494-4138-8392-4b645d3171bd.md."

# widen the "Error Detail" column (16th column) to fit the new message
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("H4").Value = "2016-09-07 06:38:15"
$wsDeDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a054dbbfb576e41805903624e40ef6b845f30dbd/e2e/da49f725-a87d-4da0-87bd-44b5adf3ff8e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3ffcc59c70c94baa481c8c9b824af0f0aa261ee/e2e/da49f725-a87d-4da0-87bd-44b5adf3ff8e.md."

$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("H5").Value = "2016-09-07 06:38:15"
$wsDeDe.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a054dbbfb576e41805903624e40ef6b845f30dbd/e2e/df1634af-b494-4138-8392-4b645d3171bd.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3ffcc59c70c94baa481c8c9b824af0f0aa261ee/e2e/df1634af-b494-4138-8392-4b645d3171bd.md."

# widen the "Error Detail" column (16th column) to fit the new message
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
